# feat: add initial parsing of world queries from data
#
# Update the "When" value of the initial card (E2) from a bare
# "environment" marker to a parseable world-query expression, widen the
# newly-important "Text of card" / "When" columns (D & E) so the content
# is readable, and move the active selection/view onto the edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Expand the "environment" trigger into a parsed world-query range
#    expression (cell E2, shared string used by the initial_card row).
$ws.Range("E2").Value = "environment=0-100,Init=0-0"

# 2. Widen column D (Text of card) and column E (When) so the longer
#    query strings are readable, matching the widths already used for
#    columns H/I.
$ws.Columns.Item(4).ColumnWidth = 33
$ws.Columns.Item(5).ColumnWidth = 44.833333333333336

# 3. Move the view back to the top-left of the sheet and select the
#    cell that was just edited.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E7").Select() | Out-Null
